# Generate Report for Handoff
# Adds two new localization files (37cbfd96... and dc0f7e0a...) with status
# "Ready for handoff" to the Overview / zh-cn / de-de sheets, pushing the
# existing ".localization-config" row down.

$wb = $excel.ActiveWorkbook

$commit = "4f09da6566c243c75ee63c1185c3815b126c78aa"
$zhCommit = "8b56d8a6e828b9e605caaab38560f23dfbd30ab6"
$deCommit = "020d89bcaa63444d4f1c5d760a7e80920619c7fc"

$hyperlinkColor = 15570276  # BGR for FF6495ED (matches existing "HyperLink" style)

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDateText($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Wipe all existing hyperlinks (Range.Hyperlinks.Delete clears the whole
# sheet's collection) so we can rebuild them cleanly in the right order.
$ws1.Range("A1").Hyperlinks.Delete()

# Shift the ".localization-config" row from row 4 down to row 6, and write
# the two new rows at 4 and 5.
$ws1.Range("A6").Value = ".localization-config"
$ws1.Range("B6").Value = "Not to be localized"
$ws1.Range("C6").Value = "Not to be localized"
Style-AsHyperlink($ws1.Range("A6"))

$ws1.Range("A4").Value = "37cbfd96-0b24-48d9-8a23-abd3b10c8770.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
Style-AsHyperlink($ws1.Range("A4"))

$ws1.Range("A5").Value = "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
Style-AsHyperlink($ws1.Range("A5"))

# Rebuild hyperlinks in document order.
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/3212c038-acf4-42be-815d-fe6e6aa7e683.md", "", "", "3212c038-acf4-42be-815d-fe6e6aa7e683.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/59731b52-8d84-406e-932c-1b6764f3c0cd.md", "", "", "59731b52-8d84-406e-932c-1b6764f3c0cd.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/37cbfd96-0b24-48d9-8a23-abd3b10c8770.md", "", "", "37cbfd96-0b24-48d9-8a23-abd3b10c8770.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md", "", "", "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

# Move the ".localization-config" row from row 4 to row 6.
$ws2.Range("A6").Value = ".localization-config"
$ws2.Range("B6").Value = "Not to be localized"
$ws2.Range("D6").Value = "0001-01-01 00:00:00"
$ws2.Range("G6").Value = "0001-01-01 00:00:00"
$ws2.Range("H6").Value = "Ignored"
Style-AsHyperlink($ws2.Range("A6"))
Style-AsDateText($ws2.Range("D6"))

# New row 4: 37cbfd96...
$ws2.Range("A4").Value = "37cbfd96-0b24-48d9-8a23-abd3b10c8770.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "37cbfd96-0b24-48d9-8a23-abd3b10c8770.7f24e5be11b4f82d9719896cc60f87e888964d01.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-08 04:51:47"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"
Style-AsHyperlink($ws2.Range("A4"))
Style-AsHyperlink($ws2.Range("C4"))
Style-AsDateText($ws2.Range("D4"))

# New row 5: dc0f7e0a...
$ws2.Range("A5").Value = "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.2c3e9ec6ee1bf76cb9ec7f73e835bd43547b4d64.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-03-08 04:51:47"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"
Style-AsHyperlink($ws2.Range("A5"))
Style-AsHyperlink($ws2.Range("C5"))
Style-AsDateText($ws2.Range("D5"))

# Rebuild hyperlinks in document order for zh-cn.
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/3212c038-acf4-42be-815d-fe6e6aa7e683.md", "", "", "3212c038-acf4-42be-815d-fe6e6aa7e683.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/3212c038-acf4-42be-815d-fe6e6aa7e683.05e25289f06c99d72e0d9d12a0d0adfec28d0880.zh-cn.xlf", "", "", "3212c038-acf4-42be-815d-fe6e6aa7e683.05e25289f06c99d72e0d9d12a0d0adfec28d0880.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/59731b52-8d84-406e-932c-1b6764f3c0cd.md", "", "", "59731b52-8d84-406e-932c-1b6764f3c0cd.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/59731b52-8d84-406e-932c-1b6764f3c0cd.fd891f01dc588c11eced54ef22601e946fca375d.zh-cn.xlf", "", "", "59731b52-8d84-406e-932c-1b6764f3c0cd.fd891f01dc588c11eced54ef22601e946fca375d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/37cbfd96-0b24-48d9-8a23-abd3b10c8770.md", "", "", "37cbfd96-0b24-48d9-8a23-abd3b10c8770.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/37cbfd96-0b24-48d9-8a23-abd3b10c8770.7f24e5be11b4f82d9719896cc60f87e888964d01.zh-cn.xlf", "", "", "37cbfd96-0b24-48d9-8a23-abd3b10c8770.7f24e5be11b4f82d9719896cc60f87e888964d01.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md", "", "", "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.2c3e9ec6ee1bf76cb9ec7f73e835bd43547b4d64.zh-cn.xlf", "", "", "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.2c3e9ec6ee1bf76cb9ec7f73e835bd43547b4d64.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

# Move the ".localization-config" row from row 4 to row 6.
$ws3.Range("A6").Value = ".localization-config"
$ws3.Range("B6").Value = "Not to be localized"
$ws3.Range("D6").Value = "0001-01-01 00:00:00"
$ws3.Range("G6").Value = "0001-01-01 00:00:00"
$ws3.Range("H6").Value = "Ignored"
Style-AsHyperlink($ws3.Range("A6"))
Style-AsDateText($ws3.Range("D6"))

# New row 4: 37cbfd96...
$ws3.Range("A4").Value = "37cbfd96-0b24-48d9-8a23-abd3b10c8770.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "37cbfd96-0b24-48d9-8a23-abd3b10c8770.7f24e5be11b4f82d9719896cc60f87e888964d01.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-08 04:51:56"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"
Style-AsHyperlink($ws3.Range("A4"))
Style-AsHyperlink($ws3.Range("C4"))
Style-AsDateText($ws3.Range("D4"))

# New row 5: dc0f7e0a...
$ws3.Range("A5").Value = "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.2c3e9ec6ee1bf76cb9ec7f73e835bd43547b4d64.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-08 04:51:56"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"
Style-AsHyperlink($ws3.Range("A5"))
Style-AsHyperlink($ws3.Range("C5"))
Style-AsDateText($ws3.Range("D5"))

# Rebuild hyperlinks in document order for de-de.
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/3212c038-acf4-42be-815d-fe6e6aa7e683.md", "", "", "3212c038-acf4-42be-815d-fe6e6aa7e683.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/3212c038-acf4-42be-815d-fe6e6aa7e683.05e25289f06c99d72e0d9d12a0d0adfec28d0880.de-de.xlf", "", "", "3212c038-acf4-42be-815d-fe6e6aa7e683.05e25289f06c99d72e0d9d12a0d0adfec28d0880.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/59731b52-8d84-406e-932c-1b6764f3c0cd.md", "", "", "59731b52-8d84-406e-932c-1b6764f3c0cd.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/59731b52-8d84-406e-932c-1b6764f3c0cd.fd891f01dc588c11eced54ef22601e946fca375d.de-de.xlf", "", "", "59731b52-8d84-406e-932c-1b6764f3c0cd.fd891f01dc588c11eced54ef22601e946fca375d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/37cbfd96-0b24-48d9-8a23-abd3b10c8770.md", "", "", "37cbfd96-0b24-48d9-8a23-abd3b10c8770.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/37cbfd96-0b24-48d9-8a23-abd3b10c8770.7f24e5be11b4f82d9719896cc60f87e888964d01.de-de.xlf", "", "", "37cbfd96-0b24-48d9-8a23-abd3b10c8770.7f24e5be11b4f82d9719896cc60f87e888964d01.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md", "", "", "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.2c3e9ec6ee1bf76cb9ec7f73e835bd43547b4d64.de-de.xlf", "", "", "dc0f7e0a-765e-4e1f-a62c-f446f7a04f66.2c3e9ec6ee1bf76cb9ec7f73e835bd43547b4d64.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit/.localization-config", "", "", ".localization-config")

Write-Output "done"
